$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.019.27'
$ws.Range("E2").Value = '  +0.37%  '

$ws.Range("D3").Value = '1.635.43'
$ws.Range("E3").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.75'
$ws.Range("E5").Value = '  -0.35%  '

$ws.Range("E6").Value = '  -0.31%  '

$ws.Range("E7").Value = '  +0.51%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.252'
$ws.Range("E8").Value = '  -1.86%  '

$ws.Range("E9").Value = '  -1.55%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.67'
$ws.Range("E10").Value = '  -4.39%  '

$ws.Range("E11").Value = '  +0.10%  '

$ws.Range("D12").Value = '1.705.60'
$ws.Range("E12").Value = '  +2.60%  '

$ws.Range("D13").Value = '1.864.59'
$ws.Range("E13").Value = '  +0.05%  '

$ws.Range("E14").Value = '  -1.56%  '

$ws.Range("E15").Value = '  -2.12%  '

$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.20'
$ws.Range("E16").Value = '  -1.04%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.0₃0747'
$ws.Range("E17").Value = '  -2.08%  '

$ws.Range("D18").Value = '26.038.18'
$ws.Range("E18").Value = '  +0.42%  '

$ws.Range("E19").Value = '  +0.42%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '190.92'
$ws.Range("E20").Value = '  -0.97%  '

$ws.Range("E21").Value = '  -1.88%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.62'
$ws.Range("E22").Value = '  -2.94%  '

$ws.Range("E23").Value = '  -1.82%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.132'
$ws.Range("E24").Value = '  +0.63%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.69'
$ws.Range("E25").Value = '  +0.14%  '

$ws.Range("B26").Value = 'BinanceUSD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.01'
$ws.Range("E26").Value = '  +0.59%  '

$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.77'
$ws.Range("E27").Value = '  -1.43%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.77'
$ws.Range("E28").Value = '  -1.68%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.30'
$ws.Range("E29").Value = '  -1.60%  '

$ws.Range("E30").Value = '  -0.46%  '

$ws.Range("E31").Value = '  -2.99%  '

$ws.Range("E32").Value = '  -2.32%  '

$ws.Range("E33").Value = '  -3.03%  '

$ws.Range("E34").Value = '  -1.56%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.880'
$ws.Range("E36").Value = '  -2.20%  '

$ws.Range("D37").Value = '1.130.67'
$ws.Range("E37").Value = '  -0.33%  '

$ws.Range("E38").Value = '  -0.21%  '

$ws.Range("E39").Value = '  -2.92%  '

$ws.Range("E40").Value = '  -0.82%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '99.01'
$ws.Range("E41").Value = '  -0.47%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.793'
$ws.Range("E42").Value = '  -0.57%  '

$ws.Range("E43").Value = '  -3.04%  '

$ws.Range("E44").Value = '  -0.88%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '55.52'
$ws.Range("E45").Value = '  -1.82%  '

$ws.Range("E46").Value = '  -0.71%  '

$ws.Range("E47").Value = '  +1.76%  '

$ws.Range("E48").Value = '  +0.04%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.56'
$ws.Range("E49").Value = '  -1.42%  '

$ws.Range("E50").Value = '  +0.33%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0929'
$ws.Range("E51").Value = '  -3.04%  '
